$wb = $excel.ActiveWorkbook

# --- Sheet 1: LH_Review_WF_LOGIN -----------------------------------------
$ws1 = $wb.Worksheets.Item("LH_Review_WF_LOGIN")

# Reviewer verification column (I2:I4): "closed" -> "open"
$ws1.Range("I2").Value = "open"
$ws1.Range("I3").Value = "open"
$ws1.Range("I4").Value = "open"

# --- Sheet 2: VERSION-HISTORY --------------------------------------------
$ws2 = $wb.Worksheets.Item("VERSION-HISTORY")

# Clear out the "v1.1 / closed comments" history row (row 3), turning it
# back into a blank template row like row 4.
$ws2.Range("A3:D3").ClearContents()

$grayFill = 16382198  # RGB(F6,F8,F9) packed as BGR for the Color property

$a3 = $ws2.Range("A3")
$a3.Interior.Color = $grayFill
$a3.HorizontalAlignment = -4108
$a3.VerticalAlignment = -4108
$a3.WrapText = $true

$b3 = $ws2.Range("B3")
$b3.Interior.Color = 16777215
$b3.VerticalAlignment = -4108
$b3.WrapText = $true

$c3 = $ws2.Range("C3")
$c3.Interior.Color = $grayFill
$c3.HorizontalAlignment = 1
$c3.VerticalAlignment = -4108
$c3.WrapText = $true

$d3 = $ws2.Range("D3")
$d3.Interior.Color = $grayFill
$d3.HorizontalAlignment = 1
$d3.VerticalAlignment = -4108
$d3.WrapText = $true
$d3.NumberFormat = "d-mmm-yy"

# --- Selections / active sheet -------------------------------------------
# Before: VERSION-HISTORY tab was selected (cell G7); after: the first
# sheet (LH_Review_WF_LOGIN) is the active tab, with I11 selected, and
# VERSION-HISTORY's own remembered selection becomes C19.
$ws2.Range("C19").Select()
$ws1.Activate()
$ws1.Range("I11").Select()
